# Applies numeric corrections to Sheets/Leviathan_Profits.xlsx (scheduled-runner update).
# Each assignment below mirrors one cell delta from the commit diff, grouped by worksheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 391.2
$ws.Range("I11").Value = 391.2
$ws.Range("K11").Value = 391.2
$ws.Range("M11").Value = -251.2
$ws.Range("H51").Value = 11908691
$ws.Range("I51").Value = 3999.5
$ws.Range("J51").Value = 16670567
$ws.Range("K51").Value = 3999.5
$ws.Range("L51").Value = 16670567
$ws.Range("M51").Value = -3515.5
$ws.Range("N51").Value = -16671535
$ws.Range("H62").Value = 8155.4
$ws.Range("I62").Value = 5099.8
$ws.Range("K62").Value = 5099.8
$ws.Range("M62").Value = -4475.8
$ws.Range("H65").Value = 8155.4
$ws.Range("I65").Value = 5099.8
$ws.Range("K65").Value = 25499
$ws.Range("M65").Value = -22379
$ws.Range("H80").Value = 1649.7307
$ws.Range("I80").Value = 1007.55554
$ws.Range("J80").Value = 1989.7059
$ws.Range("K80").Value = 3022.66662
$ws.Range("L80").Value = 5969.1177
$ws.Range("M80").Value = -2024.66662
$ws.Range("N80").Value = -7965.1177
$ws.Range("H83").Value = 1649.7307
$ws.Range("I83").Value = 1007.55554
$ws.Range("J83").Value = 1989.7059
$ws.Range("K83").Value = 9067.99986
$ws.Range("L83").Value = 17907.3531
$ws.Range("M83").Value = -4075.99986
$ws.Range("N83").Value = -27891.3531
$ws.Range("H113").Value = 51631.953
$ws.Range("I113").Value = 102376.2
$ws.Range("J113").Value = 5500.8184
$ws.Range("K113").Value = 102376.2
$ws.Range("L113").Value = 5500.8184
$ws.Range("M113").Value = -99122.2
$ws.Range("N113").Value = -12008.8184
$ws.Range("H116").Value = 10075.833
$ws.Range("I116").Value = 15209.223
$ws.Range("K116").Value = 15209.223
$ws.Range("M116").Value = -11767.223

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H26").Value = 2900
$ws.Range("I26").Value = 2900
$ws.Range("K26").Value = 2900
$ws.Range("M26").Value = -2570
$ws.Range("H32").Value = 7931
$ws.Range("I32").Value = 3052.0408
$ws.Range("K32").Value = 3052.0408
$ws.Range("M32").Value = -2765.0408
$ws.Range("H122").Value = 1444.9474
$ws.Range("I122").Value = 1382
$ws.Range("J122").Value = 2578
$ws.Range("K122").Value = 4146
$ws.Range("L122").Value = 7734
$ws.Range("M122").Value = -1696
$ws.Range("N122").Value = -12634

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1875.4706
$ws.Range("I134").Value = 1635.3636
$ws.Range("J134").Value = 2315.6667
$ws.Range("K134").Value = 4906.0908
$ws.Range("L134").Value = 6947.000100000001
$ws.Range("M134").Value = -2371.0908
$ws.Range("N134").Value = -12017.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 90914620
$ws.Range("I7").Value = 142861470
$ws.Range("J7").Value = 7650
$ws.Range("K7").Value = 142861470
$ws.Range("L7").Value = 7650
$ws.Range("M7").Value = -142861357
$ws.Range("N7").Value = -7876
$ws.Range("H22").Value = 666.4286
$ws.Range("I22").Value = 356.25
$ws.Range("J22").Value = 1080
$ws.Range("K22").Value = 356.25
$ws.Range("L22").Value = 1080
$ws.Range("M22").Value = -6.25
$ws.Range("N22").Value = -1780
$ws.Range("H31").Value = 12317.526
$ws.Range("I31").Value = 3423.3914
$ws.Range("K31").Value = 3423.3914
$ws.Range("M31").Value = -3128.3914
$ws.Range("H34").Value = 12317.526
$ws.Range("I34").Value = 3423.3914
$ws.Range("K34").Value = 3423.3914
$ws.Range("M34").Value = -3221.3914
$ws.Range("H51").Value = 14999.444
$ws.Range("J51").Value = 14999.444
$ws.Range("L51").Value = 14999.444
$ws.Range("N51").Value = -16471.444
$ws.Range("H61").Value = 14999.444
$ws.Range("J61").Value = 14999.444
$ws.Range("L61").Value = 14999.444
$ws.Range("N61").Value = -15695.444
$ws.Range("H62").Value = 2999.75
$ws.Range("I62").Value = 3228.2856
$ws.Range("K62").Value = 3228.2856
$ws.Range("M62").Value = -2604.2856
$ws.Range("H65").Value = 2999.75
$ws.Range("I65").Value = 3228.2856
$ws.Range("K65").Value = 16141.428
$ws.Range("M65").Value = -13021.428
$ws.Range("H99").Value = 16380.19
$ws.Range("I99").Value = 21078.666
$ws.Range("K99").Value = 21078.666
$ws.Range("M99").Value = -19580.666
$ws.Range("H107").Value = 3110.6667
$ws.Range("I107").Value = 3209.8572
$ws.Range("J107").Value = 2971.8
$ws.Range("K107").Value = 3209.8572
$ws.Range("L107").Value = 2971.8
$ws.Range("M107").Value = -1289.8572
$ws.Range("N107").Value = -6811.8
$ws.Range("I122").Value = 125945.125
$ws.Range("J122").Value = 8795.4
$ws.Range("K122").Value = 377835.375
$ws.Range("L122").Value = 26386.2
$ws.Range("M122").Value = -375385.375
$ws.Range("N122").Value = -31286.2
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -54920
$ws.Range("H126").Value = 16380.19
$ws.Range("I126").Value = 21078.666
$ws.Range("K126").Value = 63235.99800000001
$ws.Range("M126").Value = -60765.99800000001
$ws.Range("H132").Value = 3020.3713
$ws.Range("I132").Value = 2930.1936
$ws.Range("J132").Value = 3719.25
$ws.Range("K132").Value = 8790.5808
$ws.Range("L132").Value = 11157.75
$ws.Range("M132").Value = -6260.5808
$ws.Range("N132").Value = -16217.75
$ws.Range("H134").Value = 2790.4
$ws.Range("I134").Value = 1842
$ws.Range("J134").Value = 5635.6
$ws.Range("K134").Value = 5526
$ws.Range("L134").Value = 16906.8
$ws.Range("M134").Value = -2991
$ws.Range("N134").Value = -21976.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1545.7858
$ws.Range("I5").Value = 800.3
$ws.Range("J5").Value = 3409.5
$ws.Range("K5").Value = 2400.9
$ws.Range("L5").Value = 10228.5
$ws.Range("M5").Value = -2288.9
$ws.Range("N5").Value = -10452.5
$ws.Range("H92").Value = 1059.4286
$ws.Range("J92").Value = 452
$ws.Range("L92").Value = 1356
$ws.Range("N92").Value = -3852
$ws.Range("H131").Value = 2173.125
$ws.Range("I131").Value = 3015
$ws.Range("J131").Value = 2052.8572
$ws.Range("K131").Value = 9045
$ws.Range("L131").Value = 6158.571599999999
$ws.Range("M131").Value = -4005
$ws.Range("N131").Value = -16238.5716
$ws.Range("H135").Value = 1545.7858
$ws.Range("I135").Value = 800.3
$ws.Range("J135").Value = 3409.5
$ws.Range("K135").Value = 7202.7
$ws.Range("L135").Value = 30685.5
$ws.Range("M135").Value = -4667.7
$ws.Range("N135").Value = -35755.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2102.25
$ws.Range("I126").Value = 2004
$ws.Range("K126").Value = 6012
$ws.Range("M126").Value = -3542

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 46477.637
$ws.Range("I16").Value = 969.125
$ws.Range("J16").Value = 167833.67
$ws.Range("K16").Value = 969.125
$ws.Range("L16").Value = 167833.67
$ws.Range("M16").Value = -799.125
$ws.Range("N16").Value = -168173.67
$ws.Range("H55").Value = 431.05554
$ws.Range("I55").Value = 465.66666
$ws.Range("K55").Value = 465.66666
$ws.Range("M55").Value = -292.66666
$ws.Range("H136").Value = 3695.0322
$ws.Range("I136").Value = 3355.3635
$ws.Range("J136").Value = 4525.3335
$ws.Range("K136").Value = 10066.0905
$ws.Range("L136").Value = 13576.0005
$ws.Range("M136").Value = -7516.0905
$ws.Range("N136").Value = -18676.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 836.43335
$ws.Range("I122").Value = 772.7778
$ws.Range("K122").Value = 2318.3334
$ws.Range("M122").Value = 131.6666
$ws.Range("H136").Value = 3043.2593
$ws.Range("I136").Value = 2548.6667
$ws.Range("K136").Value = 7646.000100000001
$ws.Range("M136").Value = -5096.000100000001
